# Updates cryptos list prices / 1h volume percentages (GitHub Actions style refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value that must be kept as TEXT even when it looks numeric
# (prices like "233.40" would otherwise be auto-converted to a Double and lose
# the trailing zero / intended text formatting), then restore the cell's
# original "Normal" style so no stray formatting is left behind.
function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "37.316.43"
$ws.Range("E2").Value = "  -0.29%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "2.062.85"
$ws.Range("E3").Value = "  -0.15%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.04%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "233.40"
$ws.Range("E5").Value = "  -1.19%  "

# Row 6 - XRP
Set-TextValue $ws.Range("D6") "0.622"
$ws.Range("E6").Value = "  +1.06%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.00%  "

# Row 8 - Solana
Set-TextValue $ws.Range("D8") "56.75"
$ws.Range("E8").Value = "  -1.69%  "

# Row 9 - Cardano
Set-TextValue $ws.Range("D9") "0.383"
$ws.Range("E9").Value = "  +0.29%  "

# Row 10 - Dogecoin
Set-TextValue $ws.Range("D10") "0.0764"
$ws.Range("E10").Value = "  +0.44%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.81%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D12") "2.367.50"
$ws.Range("E12").Value = "  +0.04%  "

# Row 13 - Chainlink
Set-TextValue $ws.Range("D13") "14.41"
$ws.Range("E13").Value = "  +0.74%  "

# Row 14 - Avalanche
Set-TextValue $ws.Range("D14") "20.67"
$ws.Range("E14").Value = "  -0.87%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  -0.29%  "

# Row 16 - Polkadot
Set-TextValue $ws.Range("D16") "5.15"
$ws.Range("E16").Value = "  -0.70%  "

# Row 17 - WrappedEther
Set-TextValue $ws.Range("D17") "2.064.27"
$ws.Range("E17").Value = "  -0.06%  "

# Row 18 - WrappedBTC
Set-TextValue $ws.Range("D18") "37.274.28"
$ws.Range("E18").Value = "  -0.74%  "

# Row 19 - Uniswap
Set-TextValue $ws.Range("D19") "6.36"
$ws.Range("E19").Value = "  +3.61%  "

# Row 20 - Litecoin
Set-TextValue $ws.Range("D20") "69.48"
$ws.Range("E20").Value = "  +1.42%  "

# Row 21 - ShibaInu
$ws.Range("E21").Value = "  +0.33%  "

# Row 22 - BitcoinCash
Set-TextValue $ws.Range("D22") "226.13"
$ws.Range("E22").Value = "  +0.48%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  +1.09%  "

# Row 25 - PancakeSwap
$ws.Range("E25").Value = "  -2.48%  "

# Row 26 - Monero
Set-TextValue $ws.Range("D26") "166.25"
$ws.Range("E26").Value = "  +1.69%  "

# Row 27 - Cosmos
$ws.Range("E27").Value = "  -0.78%  "

# Row 28 - ImmutableX
$ws.Range("E28").Value = "  +2.55%  "

# Row 29 - EthereumClassic
Set-TextValue $ws.Range("D29") "19.02"
$ws.Range("E29").Value = "  -0.88%  "

# Row 30 - Kaspa
$ws.Range("E30").Value = "  -2.79%  "

# Row 31 - Stellar
Set-TextValue $ws.Range("D31") "0.118"
$ws.Range("E31").Value = "  -0.41%  "

# Row 32 - Filecoin
Set-TextValue $ws.Range("D32") "4.50"
$ws.Range("E32").Value = "  +0.84%  "

# Row 33 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D33") "4.61"
$ws.Range("E33").Value = "  +3.42%  "

# Row 34 - Hedera
Set-TextValue $ws.Range("D34") "0.0617"
$ws.Range("E34").Value = "  -1.54%  "

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = "  -5.79%  "

# Row 37 - WEMIXToken
Set-TextValue $ws.Range("D37") "1.78"
$ws.Range("E37").Value = "  -0.87%  "

# Row 38 - RenderToken
$ws.Range("E38").Value = "  -4.42%  "

# Row 39 - THORChain
Set-TextValue $ws.Range("D39") "5.56"
$ws.Range("E39").Value = "  -4.82%  "

# Row 40 - HuobiToken
$ws.Range("E40").Value = "  -0.95%  "

# Row 41 - Maker
Set-TextValue $ws.Range("D41") "1.470.52"
$ws.Range("E41").Value = "  +0.25%  "

# Row 42 - Aave
Set-TextValue $ws.Range("D42") "95.93"
$ws.Range("E42").Value = "  +0.54%  "

# Row 43 - was TrustWalletToken, now Cronos (rows 43/44 swapped identity)
$ws.Range("B43").Value = "Cronos"
$ws.Range("C43").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D43") "0.0932"
$ws.Range("E43").Value = "  -2.71%  "

# Row 44 - was Cronos, now TrustWalletToken
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D44") "1.17"
$ws.Range("E44").Value = "  +3.27%  "

# Row 45 - FTXToken
Set-TextValue $ws.Range("D45") "4.27"
$ws.Range("E45").Value = "  -4.20%  "

# Row 46 - VeChain
$ws.Range("E46").Value = "  +0.23%  "

# Row 47 - ARBITRUM
$ws.Range("E47").Value = "  -0.25%  "

# Row 48 - InjectiveProtocol
Set-TextValue $ws.Range("D48") "15.14"
$ws.Range("E48").Value = "  -5.92%  "

# Row 49 - FraxShare
$ws.Range("E49").Value = "  -1.41%  "

# Row 50 - MXToken
$ws.Range("E50").Value = "  +1.32%  "

# Row 51 - RocketPoolETH
Set-TextValue $ws.Range("D51") "2.257.38"
$ws.Range("E51").Value = "  +0.10%  "
